$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a single BV range in columns B:F (BV1, BV2, BV20,
# BV21, BV22). We need to introduce a new BV40 column as the new column B,
# pushing the existing BV1..BV22 range one column to the right (C:G), and
# append two new device rows (6 and 7).

# Step 1: shift the existing B1:F5 block one column to the right (C1:G5),
# carrying over both values and formatting (including the still-blank
# cells in row 5) so nothing gets lost.
$ws.Range("B1:F5").Copy($ws.Range("C1:G5"))

# Step 2: populate the new column B with the BV40 data, matching the style
# of the header row that now lives in C1.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "BV40"
$ws.Range("B2").Value = "inactive"
$ws.Range("B3").Value = "inactive"
$ws.Range("B4").Value = "inactive"
# B5 is left untouched - it was already blank before the shift, so it stays blank.

# Step 3: add the new deviceInstance 1601 row.
$ws.Range("A6").Value = 1601
$ws.Range("B6").Value = "inactive"
$ws.Range("C6").Font.Bold = $false
$ws.Range("D6").Font.Bold = $false
$ws.Range("E6").Value = "inactive"
$ws.Range("F6").Value = "inactive"
$ws.Range("G6").Value = "active"

# Step 4: add the new deviceInstance 1603 row.
$ws.Range("A7").Value = 1603
$ws.Range("B7").Value = "inactive"
$ws.Range("C7").Value = "inactive"
$ws.Range("D7").Value = "inactive"
$ws.Range("E7").Value = "inactive"
$ws.Range("F7").Value = "inactive"
$ws.Range("G7").Value = "inactive"
